# Rename the worksheet from "Sheet1" to "Shivam Mavi"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Shivam Mavi"

# Insert a new column at A, shifting the existing teamName..result columns
# from A..L to B..M
$ws.Columns.Item(1).Insert()

# New header for the inserted column
$ws.Range("A1").Value = "matchNo"

# matchNo value for the existing (now row 2) record
$ws.Range("A2").Value = "Final"

# New row 3 - a second match record for Shivam Mavi
$ws.Range("A3").Value = "18th"
$ws.Range("B3").Value = "Kolkata Knight Riders"
$ws.Range("C3").Value = "Shivam Mavi"
$ws.Range("D3").Value = "b Morris"
$ws.Range("E3").Value = "'5"
$ws.Range("F3").Value = "'7"
$ws.Range("G3").Value = "'1"
$ws.Range("H3").Value = "'0"
$ws.Range("I3").Value = "'71.42"
$ws.Range("J3").Value = "Rajasthan Royals"
$ws.Range("K3").Value = "Wankhede"
$ws.Range("L3").Value = "April 24"
$ws.Range("M3").Value = "Royals won by 6 wickets (with 7 balls remaining)"
